# Graph try editor excel reader added
# Inserts a new "GraphPageContent" worksheet (right after "StackPageContent")
# containing the topic/code_type/expected_result scenarios used by the
# graph page's try-editor tests.

$wb = $excel.ActiveWorkbook

$stackSheet = $wb.Worksheets.Item("StackPageContent")

$graphSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $stackSheet)
$graphSheet.Name = "GraphPageContent"

$data = @(
    @("topic_page", "code_type", "expected_result"),
    @("Graph", "valid", "Hello"),
    @("Graph", "invalid", "an error popup stating NameError: name 'invalid' is not defined on line1"),
    @("Graph Representations", "valid", "Hello"),
    @("Graph Representations", "invalid", "an error popup stating NameError: name 'invalid' is not defined on line1")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $graphSheet.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

[void]$graphSheet.Columns("C:C").Select()
